$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 46, shifting rows 46:61 down to 47:62
$ws.Rows.Item(46).Insert()

$ws.Cells.Item(46, 1).Value = 4
$ws.Cells.Item(46, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(46, 3).Value = "Los Lagos"
$ws.Cells.Item(46, 4).Value = 44964
$ws.Cells.Item(46, 5).Value = 10
$ws.Cells.Item(46, 6).Value = 100112030
$ws.Cells.Item(46, 7).Value = "Poroto granado"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 80
$ws.Cells.Item(46, 11).Value = 42000
$ws.Cells.Item(46, 12).Value = 42000
$ws.Cells.Item(46, 13).Value = 42000
$ws.Cells.Item(46, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(46, 15).Value = "Región Metropolitana"
$ws.Cells.Item(46, 16).Value = 1680
$ws.Cells.Item(46, 17).Value = 25
$ws.Cells.Item(46, 18).Value = "Hortaliza"
